$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell P4 (year 2022) - copy format from O4
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Value = 2022

# Row 5 (Total row) - copy format from O5, bold/Times font style
$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("P5").Value = 1188.7
$ws.Range("P5").NumberFormat = "#,##0.0"
$ws.Range("P5").HorizontalAlignment = -4152
$ws.Range("P5").VerticalAlignment = -4107

# Row 6
$ws.Range("O6").Copy()
$ws.Range("P6").PasteSpecial(-4122)
$ws.Range("P6").Value = 263.89999999999998
$ws.Range("P6").NumberFormat = "#,##0.0"
$ws.Range("P6").HorizontalAlignment = -4152
$ws.Range("P6").VerticalAlignment = -4107

# Row 7
$ws.Range("O7").Copy()
$ws.Range("P7").PasteSpecial(-4122)
$ws.Range("P7").Value = 263.2
$ws.Range("P7").NumberFormat = "#,##0.0"
$ws.Range("P7").HorizontalAlignment = -4152
$ws.Range("P7").VerticalAlignment = -4107

# Row 8
$ws.Range("O8").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("P8").Value = 12.4
$ws.Range("P8").NumberFormat = "#,##0.0"
$ws.Range("P8").HorizontalAlignment = -4152
$ws.Range("P8").VerticalAlignment = -4107

# Row 9 (text "-")
$ws.Range("O9").Copy()
$ws.Range("P9").PasteSpecial(-4122)
$ws.Range("P9").Value = "-"
$ws.Range("P9").NumberFormat = "#,##0.0"
$ws.Range("P9").HorizontalAlignment = -4152
$ws.Range("P9").VerticalAlignment = -4107

# Row 10
$ws.Range("O10").Copy()
$ws.Range("P10").PasteSpecial(-4122)
$ws.Range("P10").Value = 93
$ws.Range("P10").NumberFormat = "#,##0.0"
$ws.Range("P10").HorizontalAlignment = -4152
$ws.Range("P10").VerticalAlignment = -4107

# Row 11
$ws.Range("O11").Copy()
$ws.Range("P11").PasteSpecial(-4122)
$ws.Range("P11").Value = 171.5
$ws.Range("P11").NumberFormat = "#,##0.0"
$ws.Range("P11").HorizontalAlignment = -4152
$ws.Range("P11").VerticalAlignment = -4107

# Row 12
$ws.Range("O12").Copy()
$ws.Range("P12").PasteSpecial(-4122)
$ws.Range("P12").Value = 220.6
$ws.Range("P12").NumberFormat = "#,##0.0"
$ws.Range("P12").HorizontalAlignment = -4152
$ws.Range("P12").VerticalAlignment = -4107

# Row 13
$ws.Range("O13").Copy()
$ws.Range("P13").PasteSpecial(-4122)
$ws.Range("P13").Value = 159.30000000000001
$ws.Range("P13").NumberFormat = "#,##0.0"
$ws.Range("P13").HorizontalAlignment = -4152
$ws.Range("P13").VerticalAlignment = -4107

# Row 14
$ws.Range("O14").Copy()
$ws.Range("P14").PasteSpecial(-4122)
$ws.Range("P14").Value = 1.7
$ws.Range("P14").NumberFormat = "#,##0.0"
$ws.Range("P14").HorizontalAlignment = -4152
$ws.Range("P14").VerticalAlignment = -4107

# Row 15 (text "-")
$ws.Range("O15").Copy()
$ws.Range("P15").PasteSpecial(-4122)
$ws.Range("P15").Value = "-"
$ws.Range("P15").NumberFormat = "#,##0.0"
$ws.Range("P15").HorizontalAlignment = -4152
$ws.Range("P15").VerticalAlignment = -4107

# Row 16 (bottom border row)
$ws.Range("O16").Copy()
$ws.Range("P16").PasteSpecial(-4122)
$ws.Range("P16").Value = 3.1
$ws.Range("P16").NumberFormat = "#,##0.0"
$ws.Range("P16").HorizontalAlignment = -4152
$ws.Range("P16").VerticalAlignment = -4107

$ws.Range("Q7").Select()
